$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4: rename inference dataset headers from InferDS1-x to InferDS2-x ---
$ws.Range("E4").Value = "InferDS2-0"
$ws.Range("F4").Value = "InferDS2-1"
$ws.Range("G4").Value = "InferDS2-2"
$ws.Range("H4").Value = "InferDS2-3"

# --- Row 7: second data-set row now also carries its own SimulId (A7) ---
$ws.Range("A7").Value = 12932

# B7/C7 already equal "Engine2"/"trainDS2" - keep them explicit for clarity
$ws.Range("B7").Value = "Engine2"
$ws.Range("C7").Value = "trainDS2"

# --- Clear the old stand-alone "infer cmd" preview block (rows 8-11) ---
$ws.Range("E8:E11").ClearContents()

# --- New preview block now lives two rows further down (rows 10-13) ---
$ws.Range("E10").Value = "zzz Infer 50 Config/50/Client.xml Config/50/InferDS2-0.xml  12932"
$ws.Range("E11").Value = "zzz Infer 50 Config/50/Client.xml Config/50/InferDS2-1.xml  12932"
$ws.Range("E12").Value = "zzz Infer 50 Config/50/Client.xml Config/50/InferDS2-2.xml  12932"
$ws.Range("E13").Value = "zzz Infer 50 Config/50/Client.xml Config/50/InferDS2-3.xml  12932"

# --- Update the sheet selection to mirror the saved view state ---
$ws.Range("E10:E13").Select()

$wb.Save()
